$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J column (situacao) text refresh: "meses sem comprar" counters advanced ---
$ws.Range("J2").Value = "INATIVO - 54.3 meses sem comprar"
$ws.Range("J4").Value = "INATIVO - 35.1 meses sem comprar"
$ws.Range("J8").Value = "INATIVO - 17.0 meses sem comprar"
$ws.Range("J10").Value = "INATIVO - 1.6 meses sem comprar"
$ws.Range("J11").Value = "INATIVO - 4.0 meses sem comprar"
$ws.Range("J17").Value = "INATIVO - 36.4 meses sem comprar"
$ws.Range("J18").Value = "INATIVO - 10.4 meses sem comprar"
$ws.Range("J19").Value = "INATIVO - 14.3 meses sem comprar"
$ws.Range("J20").Value = "INATIVO - 36.3 meses sem comprar"
$ws.Range("J25").Value = "INATIVO - 18.5 meses sem comprar"
$ws.Range("J27").Value = "INATIVO - 16.6 meses sem comprar"
$ws.Range("J30").Value = "INATIVO - 21.5 meses sem comprar"
$ws.Range("J33").Value = "INATIVO - 13.6 meses sem comprar"
$ws.Range("J38").Value = "INATIVO - 33.4 meses sem comprar"
$ws.Range("J39").Value = "INATIVO - 11.5 meses sem comprar"
$ws.Range("J45").Value = "INATIVO - 1.5 meses sem comprar"
$ws.Range("J49").Value = "INATIVO - 8.9 meses sem comprar"
$ws.Range("J57").Value = "INATIVO - 10.9 meses sem comprar"
$ws.Range("J63").Value = "INATIVO - 26.8 meses sem comprar"
$ws.Range("J64").Value = "INATIVO - 20.9 meses sem comprar"
$ws.Range("J66").Value = "INATIVO - 11.5 meses sem comprar"
$ws.Range("J67").Value = "INATIVO - 10.3 meses sem comprar"
$ws.Range("J71").Value = "INATIVO - 0.3 meses sem comprar"
$ws.Range("J72").Value = "INATIVO - 20.9 meses sem comprar"
$ws.Range("J73").Value = "INATIVO - 32.4 meses sem comprar"
$ws.Range("J74").Value = "INATIVO - 6.7 meses sem comprar"
$ws.Range("J79").Value = "INATIVO - 21.2 meses sem comprar"
$ws.Range("J81").Value = "INATIVO - 25.2 meses sem comprar"
$ws.Range("J83").Value = "INATIVO - 20.7 meses sem comprar"
$ws.Range("J84").Value = "INATIVO - 8.4 meses sem comprar"
$ws.Range("J85").Value = "INATIVO - 14.0 meses sem comprar"
$ws.Range("J86").Value = "INATIVO - 3.9 meses sem comprar"
$ws.Range("J87").Value = "INATIVO - 10.7 meses sem comprar"
$ws.Range("J88").Value = "INATIVO - 10.1 meses sem comprar"
$ws.Range("J89").Value = "INATIVO - 14.1 meses sem comprar"
$ws.Range("J90").Value = "INATIVO - 32.4 meses sem comprar"
$ws.Range("J91").Value = "INATIVO - 12.7 meses sem comprar"
$ws.Range("J92").Value = "INATIVO - 17.7 meses sem comprar"
$ws.Range("J93").Value = "INATIVO - 15.4 meses sem comprar"
$ws.Range("J94").Value = "INATIVO - 18.2 meses sem comprar"
$ws.Range("J95").Value = "INATIVO - 31.9 meses sem comprar"
$ws.Range("J97").Value = "INATIVO - 1.2 meses sem comprar"
$ws.Range("J98").Value = "INATIVO - 21.8 meses sem comprar"
$ws.Range("J99").Value = "INATIVO - 36.2 meses sem comprar"
$ws.Range("J101").Value = "INATIVO - 13.5 meses sem comprar"
$ws.Range("J102").Value = "INATIVO - 23.9 meses sem comprar"
$ws.Range("J103").Value = "INATIVO - 9.5 meses sem comprar"
$ws.Range("J104").Value = "INATIVO - 24.3 meses sem comprar"
$ws.Range("J105").Value = "INATIVO - 14.7 meses sem comprar"
$ws.Range("J106").Value = "INATIVO - 5.6 meses sem comprar"
$ws.Range("J107").Value = "INATIVO - 20.3 meses sem comprar"
$ws.Range("J108").Value = "INATIVO - 6.1 meses sem comprar"
$ws.Range("J110").Value = "INATIVO - 7.5 meses sem comprar"

# --- Row 54 (id_cliente 5986): updated probabilities / dates ---
$ws.Range("B54").Value = 0.75
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 26
$ws.Range("F54").Value = 1
$ws.Range("H54").Value = 45805.76532407408
$ws.Range("I54").Value = 45866.76532407408

# --- Row 108 (id_cliente 27514): customer went INATIVO ---
$ws.Range("I108").Value = "INATIVO"
$ws.Range("I108").NumberFormat = "dd/mm/yyyy"
$ws.Range("J108").Value = "INATIVO - 6.1 meses sem comprar"

# --- Row 111 (id_cliente 28458): updated counters / dates ---
$ws.Range("E111").Value = 14943
$ws.Range("H111").Value = 45805.68260416666
